# Refactor classification synthetic data generation
#
# The underlying Python simulation that produces this report was re-run
# (refactored), which re-randomizes the synthetic "count"/"FairAI_target"/
# "FairAI_pred" columns (F/G/H) for the "equal_distribution" and
# "equal_distribution_..._poor_model" simulation blocks, and fixes a typo
# in one of the recorded target_bias labels (D5/D17/D29).
#
# (Rows 14-25, the "unequal_distribution" block, are numerically unchanged
# in the source data - only their string precision changed on save - so
# they are intentionally left untouched here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mislabeled target_bias text: "[0.34, 0.33, 0.33]" -> "[0.33, 0.33, 0.33]"
# (affects every cell sharing that string: D5, D17, D29)
$ws.Range("D5").Value  = "[0.33, 0.33, 0.33]"
$ws.Range("D17").Value = "[0.33, 0.33, 0.33]"
$ws.Range("D29").Value = "[0.33, 0.33, 0.33]"

# --- high_bias_single_cat_equal_distribution_classification (rows 2-13)
$ws.Range("F2").Value  = 1983
$ws.Range("G2").Value  = 0.9730663837782545
$ws.Range("H2").Value  = 0.9953757150787123

$ws.Range("G3").Value  = 0.9792545727230993
$ws.Range("H3").Value  = 0.9908263448068602

$ws.Range("F4").Value  = 2003
$ws.Range("G4").Value  = 0.9821234791256194
$ws.Range("H4").Value  = 0.9916921031082567

$ws.Range("F5").Value  = 2008
$ws.Range("G5").Value  = 0.9799819065811723
$ws.Range("H5").Value  = 0.8763972190758041

$ws.Range("F6").Value  = 1972
$ws.Range("G6").Value  = 0.9820925564120693
$ws.Range("H6").Value  = 0.9315771213971266

$ws.Range("F7").Value  = 2020
$ws.Range("G7").Value  = 0.9957869989718151
$ws.Range("H7").Value  = 0.9414718933330611

$ws.Range("F8").Value  = 1897
$ws.Range("G8").Value  = 0.8569894699374928
$ws.Range("H8").Value  = 0.9061259020888963

$ws.Range("F9").Value  = 2101
$ws.Range("G9").Value  = 0.9337302346033731
$ws.Range("H9").Value  = 0.9476833479013603

$ws.Range("F10").Value = 2002
$ws.Range("G10").Value = 0.9260077671576913
$ws.Range("H10").Value = 0.9579843375632696

$ws.Range("F11").Value = 1929
$ws.Range("G11").Value = 0.8742088238668536
$ws.Range("H11").Value = 0.774184365108198

$ws.Range("F12").Value = 2081
$ws.Range("G12").Value = 0.9374814945662673
$ws.Range("H12").Value = 0.882061809996029

$ws.Range("F13").Value = 1990
$ws.Range("G13").Value = 0.9383284644306384
$ws.Range("H13").Value = 0.8805044287834728

# --- high_bias_single_cat_equal_distribution_classification_poor_model (rows 26-37)
$ws.Range("F26").Value = 1983
$ws.Range("G26").Value = 0.9822400994231908
$ws.Range("H26").Value = 0.9873068028819293

$ws.Range("G27").Value = 0.9875209121486256
$ws.Range("H27").Value = 0.9848426074188429

$ws.Range("F28").Value = 2003
$ws.Range("G28").Value = 0.9708017405676321
$ws.Range("H28").Value = 0.9893454395628996

$ws.Range("F29").Value = 2008
$ws.Range("G29").Value = 0.9927928563867169
$ws.Range("H29").Value = 0.8314686377483625

$ws.Range("F30").Value = 1972
$ws.Range("G30").Value = 0.9806781541503131
$ws.Range("H30").Value = 0.9163805059369896

$ws.Range("F31").Value = 2020
$ws.Range("G31").Value = 0.9871933421244188
$ws.Range("H31").Value = 0.9109364091862266

$ws.Range("F32").Value = 1897
$ws.Range("G32").Value = 0.8733483936198283
$ws.Range("H32").Value = 0.893018296712767

$ws.Range("F33").Value = 2101
$ws.Range("G33").Value = 0.9520549035942528
$ws.Range("H33").Value = 0.9398931018317949

$ws.Range("F34").Value = 2002
$ws.Range("G34").Value = 0.9226737494192788
$ws.Range("H34").Value = 0.9532209597940299

$ws.Range("F35").Value = 1929
$ws.Range("G35").Value = 0.876566763544507
$ws.Range("H35").Value = 0.6864097780382082

$ws.Range("F36").Value = 2081
$ws.Range("G36").Value = 0.9519017535880143
$ws.Range("H36").Value = 0.8277898707885657

$ws.Range("F37").Value = 1990
$ws.Range("G37").Value = 0.925837185109776
$ws.Range("H37").Value = 0.8396855433622178

# --- Remove the conditional formatting rule on column H (H1:H1048576, "less than 0.8")
# that highlighted low-accuracy predictions; no longer needed after the refactor.
$ws.Cells.FormatConditions.Delete()
